$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns updated in each data row (E, J, O stay 0 and are untouched)
$cols = "B", "C", "D", "F", "G", "H", "I", "K", "L", "M", "N"

# New values per row, in the same column order as $cols
$rowData = @{
  2 = @(0.5715419746630062, 0.114946053896702, 0.0527335108562923, 1.152677889745412, 1.019601378747737, 1.034541761508748, 1.051614530899876, 0.3551172859349379, 0.2885973427854225, 0.1927390872979018, 2.101940947907526)
  3 = @(0.5352276124321804, 0.1116367243003751, 0.05145903917820505, 1.14961991734846, 1.018297715806369, 1.038286947416836, 1.055839669135494, 0.3185022999422245, 0.285204863958235, 0.1854537833095442, 2.121215455152647)
  4 = @(0.5131681245110542, 0.1095733237547591, 0.05066487978112022, 1.148355659305601, 1.018066036729778, 1.04100366414454, 1.0588995300234, 0.2961017084816433, 0.2832640314981276, 0.1810748712291819, 2.133681894355242)
  5 = @(0.5042389341452349, 0.1087245520908766, 0.0503383420126049, 1.1479947510811, 1.018114621816494, 1.042215733712652, 1.060263577577402, 0.2869940166240355, 0.2825089420122566, 0.1793142401083649, 2.138920936537886)
  6 = @(0.5027599003184946, 0.1085831356539089, 0.05028394520080681, 1.147944142768665, 1.018131325186829, 1.042423340401612, 1.060497152397389, 0.2854829529700567, 0.2823857251328903, 0.1790233297350987, 2.139800470491494)
  7 = @(0.513047457842589, 0.1095619089985362, 0.05066048774962439, 1.148350167195943, 1.018066113016047, 1.041019585341502, 1.058917451730565, 0.2959787945962375, 0.283253702991594, 0.1810510301989439, 2.133751906676508)
  8 = @(0.5589717188142345, 0.1138115203462178, 0.05229649385659485, 1.15149622936476, 1.019033806892793, 1.035746571990032, 1.052974749047543, 0.3424757768959807, 0.2873981487747699, 0.1902076011206546, 2.108455662085351)
  9 = @(0.6508995031773566, 0.1218960973054948, 0.05541199685193732, 1.162531743675686, 1.025447198962567, 1.0287128801763, 1.045014041999913, 0.4342906256321442, 0.2966511855519656, 0.2089085003230835, 2.063862604341644)
  10 = @(0.7195657528351376, 0.1276855559148373, 0.05764405259679251, 1.17360808165003, 1.032918118447895, 1.02555749382978, 1.04141556215103, 0.5021277950430374, 0.3041338423120408, 0.2230993731717348, 2.034154513700287)
  11 = @(0.7510461131251702, 0.1302870355590073, 0.05864704441080448, 1.179292211448399, 1.036917382022139, 1.024558291377105, 1.040267001787157, 0.5330705122367476, 0.3076861934732023, 0.229652678132453, 2.0213024620918)
  12 = @(0.763001586552349, 0.1312675371417384, 0.0590250602820106, 1.181537455098663, 1.038518255956774, 1.024242577201065, 1.03990228011132, 0.5447994512236392, 0.3090526727315819, 0.232148230330111, 2.01653098208196)
  13 = @(0.7604252331250052, 0.1310565740935914, 0.05894372788672086, 1.181049775419652, 1.038169633749632, 1.024307785975125, 1.039977706996694, 0.5422729041460173, 0.308757431556316, 0.2316101494325764, 2.017554365395057)
  14 = @(0.7520290076971605, 0.1303677946271904, 0.05867818007186543, 1.179475069563381, 1.037047354168664, 1.024531061897235, 1.040235588845647, 0.5340352286622476, 0.3077981884225949, 0.2298577095885364, 2.020907998911802)
  15 = @(0.7468905608612602, 0.12994529548223, 0.05851529011613366, 1.178522599545971, 1.036371184255799, 1.024675983451232, 1.040402692199024, 0.528990914770759, 0.3072133937429982, 0.2287861041414914, 2.022974611084098)
  16 = @(0.717513298733536, 0.1275148960476145, 0.05757825429988372, 1.17324959611544, 1.032668849933472, 1.025631566114953, 1.041500449618709, 0.500107259481581, 0.3039046693922387, 0.2226730582809822, 2.03500775908169)
  17 = @(0.6995533680427002, 0.126015687249847, 0.05700023135200638, 1.170180074533036, 1.030551484720149, 1.026329461228897, 1.042298967126072, 0.4824091332403952, 0.3019128493705523, 0.218947877309887, 2.042559390280417)
  18 = @(0.6892462630169405, 0.1251503533735558, 0.05666660332138207, 1.168475321201399, 1.02939016991688, 1.026771932819983, 1.042804225480431, 0.4722374968988277, 0.3007811842632009, 0.216814464058892, 2.04696522374735)
  19 = @(0.6857604152193062, 0.1248568458879049, 0.05655344304992838, 1.167908556461384, 1.029006677278161, 1.026928800251156, 1.042983193636445, 0.4687949163706548, 0.3004004245004523, 0.2160937124558444, 2.048467667139796)
  20 = @(0.7014628594582746, 0.1261755940055451, 0.05706188351326347, 1.170500542700992, 1.030771029995776, 1.026250919920386, 1.042209205760834, 0.4842923184021402, 0.3021234359170109, 0.2193434765148368, 2.041749054284953)
  21 = @(0.7544942491412314, 0.1305702312661055, 0.05875622674640368, 1.179935080852459, 1.037374648471058, 1.024463780211548, 1.040157937278572, 0.5364545216528427, 0.3080793643538158, 0.2303720656885986, 2.019920368083476)
  22 = @(0.7893544432070883, 0.1334154403783856, 0.0598531063485197, 1.186641917837093, 1.042194358306546, 1.023660995165713, 1.039226562518493, 0.570612999963771, 0.3120959187369579, 0.2376611936446125, 2.006209653227288)
  23 = @(0.7707306665348881, 0.1318993623859228, 0.05926864413976318, 1.183012876782357, 1.039575865500851, 1.024056060102339, 1.039686214734417, 0.5523759260026395, 0.309940883050885, 0.2337634444778161, 2.013476457792699)
  24 = @(0.7005995213541496, 0.1261033108123257, 0.05703401466588076, 1.170355472281116, 1.030671599184828, 1.026286299970465, 1.042249643022593, 0.4834409202342442, 0.3020281878125246, 0.2191646004309433, 2.042115206896643)
  25 = @(0.625831738662697, 0.1197355004622835, 0.05457913552403681, 1.159025183985577, 1.023228136179895, 1.030262049973004, 1.046772382684658, 0.4093849714951716, 0.2940276623286024, 0.2037699531850805, 2.075389720338695)
}

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $r).Value = $vals[$i]
    }
}

Write-Output "applied updates to rows 2-25"